$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item(1)
$wsCarbon = $wb.Worksheets.Item(2)

# --- Rename "Sheet1" -> "carbon_intensity" ---
$wsCarbon.Name = "carbon_intensity"

# --- params sheet: clear the second "carbon_intensity" scenario row (row 7),
#     keeping I7's date-style formatting but no value. ---
$wsParams.Range("A7:J7").ClearContents()
$wsParams.Range("S7").ClearContents()

# --- carbon_intensity sheet: insert a "scenario" column and a "ref value" rename,
#     append an "id" column header, shifting existing data one column to the right. ---

# shift "variability growth" column (E) -> F
$wsCarbon.Range("F1").Value2 = $wsCarbon.Range("E1").Value2
$wsCarbon.Range("F2").Value2 = $wsCarbon.Range("E2").Value2
$wsCarbon.Range("F3").Value2 = $wsCarbon.Range("E3").Value2

# shift "initial_value_proportional_variation" column (D) -> E
$wsCarbon.Range("E1").Value2 = $wsCarbon.Range("D1").Value2
$wsCarbon.Range("E2").Value2 = $wsCarbon.Range("D2").Value2
$wsCarbon.Range("E3").Value2 = $wsCarbon.Range("D3").Value2

# "mean growth" column (C) -> D, with updated growth values
$wsCarbon.Range("D1").Value2 = $wsCarbon.Range("C1").Value2
$wsCarbon.Range("D2").Value2 = 0.1
$wsCarbon.Range("D3").Value2 = 1

# old "value" column (B) -> C, renamed header "ref value"
$wsCarbon.Range("C1").Value2 = "ref value"
$wsCarbon.Range("C2").Value2 = $wsCarbon.Range("B2").Value2
$wsCarbon.Range("C3").Value2 = $wsCarbon.Range("B3").Value2

# new "scenario" column (B), left blank for data rows
$wsCarbon.Range("B1").Value2 = "scenario"
$wsCarbon.Range("B2").ClearContents()
$wsCarbon.Range("B3").ClearContents()

# old "key" column (A) renamed "region"
$wsCarbon.Range("A1").Value2 = "region"

# new trailing "id" column header
$wsCarbon.Range("G1").Value2 = "id"

# --- selections: set params!F5 first, then activate carbon_intensity!D3 last
#     so "carbon_intensity" remains the active tab. ---
$wsParams.Range("F5").Select()
$wsCarbon.Range("D3").Select()
